# The deck's theme parts (ppt/theme/theme1.xml = "Office Theme", used by the
# Notes Master, and ppt/theme/theme2.xml = "Integral", used by the
# Slide Master / overall presentation design) were swapped: the Slide
# Master's design now uses the plain Office colour palette, while the
# Integral palette moves over to the Notes Master's theme.
#
# The PowerPoint object model only exposes a writable colour scheme on the
# slide-facing theme (via Slide/SlideRange/CustomLayout .ThemeColorScheme,
# which all point at the single design theme applied to the slides/master),
# so we recolour that shared theme to the target ("Office Theme") palette,
# one swatch at a time, in clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$slides = $p.Slides.Range()
$colors = $slides.ThemeColorScheme

$colors.Item(1).RGB  = RGB(0, 0, 0)        # dk1      000000
$colors.Item(2).RGB  = RGB(255, 255, 255)  # lt1      FFFFFF
$colors.Item(3).RGB  = RGB(68, 84, 106)    # dk2      44546A
$colors.Item(4).RGB  = RGB(231, 230, 230)  # lt2      E7E6E6
$colors.Item(5).RGB  = RGB(91, 155, 213)   # accent1  5B9BD5
$colors.Item(6).RGB  = RGB(237, 125, 49)   # accent2  ED7D31
$colors.Item(7).RGB  = RGB(165, 165, 165)  # accent3  A5A5A5
$colors.Item(8).RGB  = RGB(255, 192, 0)    # accent4  FFC000
$colors.Item(9).RGB  = RGB(68, 114, 196)   # accent5  4472C4
$colors.Item(10).RGB = RGB(112, 173, 71)   # accent6  70AD47
$colors.Item(11).RGB = RGB(5, 99, 193)     # hlink    0563C1
$colors.Item(12).RGB = RGB(149, 79, 114)   # folHlink 954F72
